$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 151 (shifts existing rows 151-217 down to 152-218)
$ws.Rows.Item(151).Insert()

# Populate the new row 151 with the data from the commit
$ws.Cells.Item(151, 1).Value = 4
$ws.Cells.Item(151, 2).Value = 'Feria Lagunitas de Puerto Montt'
$ws.Cells.Item(151, 3).Value = 'Los Lagos'
$ws.Cells.Item(151, 4).Value = 44609
$ws.Cells.Item(151, 5).Value = 10
$ws.Cells.Item(151, 6).Value = 100112021
$ws.Cells.Item(151, 7).Value = 'Ají'
$ws.Cells.Item(151, 8).Value = 'Inferno'
$ws.Cells.Item(151, 9).Value = 'Primera'
$ws.Cells.Item(151, 10).Value = 50
$ws.Cells.Item(151, 11).Value = 23000
$ws.Cells.Item(151, 12).Value = 23000
$ws.Cells.Item(151, 13).Value = 23000
$ws.Cells.Item(151, 14).Value = '$/caja 15 kilos'
$ws.Cells.Item(151, 15).Value = 'Limache'
$ws.Cells.Item(151, 16).Value = 1533
$ws.Cells.Item(151, 17).Value = 15
$ws.Cells.Item(151, 18).Value = 'Hortaliza'
